# Add a new row of test-suite data to the "SafeWay" sheet, trim the
# now-unused trailing blank rows, and make "SafeWay" the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SafeWay")

# Populate row 3 with the new test suite entry.
$ws.Range("A3").Value = "ON"
$ws.Range("B3").Value = "Jobs"
$ws.Range("C3").Value = 1

# Remove the now-unused trailing blank rows (4:6), shrinking the used range.
$ws.Rows.Item(4).Resize(3, $ws.Rows.Item(4).Columns.Count).EntireRow.Delete()

# Select E8 on the SafeWay sheet, then make it the active sheet/tab.
$ws.Range("E8").Select()
$ws.Activate()
